$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParaXml($paragraph, [string]$innerXml) {
    $r = $paragraph.Range
    # Exclude the trailing paragraph-mark character so the paragraph itself
    # is kept (InsertXML on a range that includes the mark would merge the
    # paragraph into the next one).
    $full = $d.Range($r.Start, $r.End - 1)
    $full.InsertXML("<w:p $wNs>$innerXml</w:p>")
}

# --- Change 1 ---
# "Being my first health data project..." paragraph: split the tail of the
# sentence into several runs and change its wording.
$p1 = $d.Paragraphs(7)
$xml1 = '<w:r><w:t>Being my first health data project, it seems an interesting dataset to start with as it appears to be both manageable and sufficiently challenging. I hope this analysis brings impactful insights on</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> patents associated with</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> FDA-approved</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> drugs.</w:t></w:r>'
Set-ParaXml $p1 $xml1

# --- Change 2 ---
# The paragraph that used to read "The first thing to check..." now carries
# the text that used to be in the following paragraph ("The next step...
# PreAnalysis."), extended with the ".ipynb" file extension and a closing
# parenthesis.
$p2 = $d.Paragraphs(11)
$xml2 = '<w:r><w:t xml:space="preserve">The next step was to dive into the dataset. To explore the dataset we opted to use pandas, a python library very useful for data manipulation and data transformation (see </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>PreAnalysis.</w:t></w:r>' +
        '<w:r><w:t>ipynb</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t>).</w:t></w:r>'
Set-ParaXml $p2 $xml2

# The paragraph that used to hold "The next step..." is now emptied out, and
# a brand-new empty paragraph (same justified-alignment formatting) follows
# it, right before the section break.
$p3 = $d.Paragraphs(12)
$xmlEmptyPair = '<w:p ' + $wNs + '><w:pPr><w:jc w:val="both"/></w:pPr></w:p>' +
                '<w:p ' + $wNs + '><w:pPr><w:jc w:val="both"/></w:pPr></w:p>'
$r3 = $p3.Range
$full3 = $d.Range($r3.Start, $r3.End - 1)
$full3.InsertXML($xmlEmptyPair)
